# Refresh FFXIV (Ixion server) Leve crafting-profit market data.
# Scheduled runner pulls the latest Universalis market-board averages and
# re-derives currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# for the affected leves on each crafting-class sheet.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H80" = 527.9474
    "I80" = 407.16666
    "J80" = 735
    "K80" = 1221.49998
    "L80" = 2205
    "M80" = -223.4999800000001
    "N80" = -4201
    "H83" = 527.9474
    "I83" = 407.16666
    "J83" = 735
    "K83" = 3664.49994
    "L83" = 6615
    "M83" = 1327.50006
    "N83" = -16599
    "H113" = 13891478
    "I113" = 2614.2856
    "K113" = 2614.2856
    "M113" = 639.7143999999998
    "H141" = 1415.2858
    "I141" = 1015.2439
    "J141" = 3465.5
    "K141" = 3045.7317
    "L141" = 10396.5
    "M141" = 2134.2683
    "N141" = -20756.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H32" = 5418.2446
    "I32" = 3615.7214
    "J32" = 18363.637
    "K32" = 3615.7214
    "L32" = 18363.637
    "M32" = -3328.7214
    "N32" = -18937.637
    "H74" = 1013.42224
    "I74" = 955.931
    "J74" = 1117.625
    "K74" = 955.931
    "L74" = 1117.625
    "M74" = -81.93100000000004
    "N74" = -2865.625
    "H77" = 1013.42224
    "I77" = 955.931
    "J77" = 1117.625
    "K77" = 4779.655000000001
    "L77" = 5588.125
    "M77" = -411.6550000000007
    "N77" = -14324.125
    "H92" = 39270
    "J92" = 39270
    "L92" = 39270
    "N92" = -44262
    "H132" = 3033477
    "I132" = 2425.2856
    "J132" = 8337817.5
    "K132" = 7275.8568
    "L132" = 25013452.5
    "M132" = -4745.8568
    "N132" = -25018512.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H76" = 42310
    "J76" = 42310
    "L76" = 42310
    "N76" = -42940
    "H79" = 42310
    "J79" = 42310
    "L79" = 42310
    "N79" = -44494
    "H80" = 282.31818
    "I80" = 78
    "J80" = 327.72223
    "K80" = 78
    "L80" = 327.72223
    "M80" = 920
    "N80" = -2323.72223
    "H83" = 282.31818
    "I83" = 78
    "J83" = 327.72223
    "K83" = 390
    "L83" = 1638.61115
    "M83" = 4602
    "N83" = -11622.61115
    "H134" = 4003.2683
    "I134" = 4654.1377
    "J134" = 2430.3333
    "K134" = 13962.4131
    "L134" = 7290.999899999999
    "M134" = -11427.4131
    "N134" = -12360.9999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 6800.8623
    "I31" = 1901.5135
    "J31" = 15433.048
    "K31" = 1901.5135
    "L31" = 15433.048
    "M31" = -1606.5135
    "N31" = -16023.048
    "H34" = 6800.8623
    "I34" = 1901.5135
    "J34" = 15433.048
    "K34" = 1901.5135
    "L34" = 15433.048
    "M34" = -1699.5135
    "N34" = -15837.048
    "H58" = 1050.7937
    "I58" = 661.90247
    "J58" = 1775.5454
    "K58" = 661.90247
    "L58" = 1775.5454
    "M58" = -458.90247
    "N58" = -2181.5454
    "H125" = 0
    "J125" = 0
    "L125" = 0
    "H136" = 1050.7937
    "I136" = 661.90247
    "J136" = 1775.5454
    "K136" = 1985.70741
    "L136" = 5326.6362
    "M136" = 564.29259
    "N136" = -10426.6362
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N125").ClearContents()  # no HQ price available; cell removed

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H64" = 2273.1428
    "I64" = 1128
    "J64" = 3800
    "K64" = 3384
    "L64" = 11400
    "M64" = -3114
    "N64" = -11940
    "H67" = 2273.1428
    "I67" = 1128
    "J67" = 3800
    "K67" = 3384
    "L67" = 11400
    "M67" = -2448
    "N67" = -13272
    "H121" = 937.069
    "I121" = 359.66666
    "J121" = 1087.6957
    "K121" = 1078.99998
    "L121" = 3263.0871
    "M121" = 231.0000199999999
    "N121" = -5883.0871
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H70" = 7688.5
    "I70" = 8168
    "J70" = 6250
    "K70" = 8168
    "L70" = 6250
    "M70" = -7898
    "N70" = -6790
    "H73" = 7688.5
    "I73" = 8168
    "J73" = 6250
    "K73" = 8168
    "L73" = 6250
    "M73" = -7232
    "N73" = -8122
    "H102" = 1235.3
    "I102" = 965.4286
    "K102" = 965.4286
    "M102" = 656.5714
    "H122" = 2161873.5
    "I122" = 2947384.8
    "J122" = 1718.125
    "K122" = 8842154.399999999
    "L122" = 5154.375
    "M122" = -8839704.399999999
    "N122" = -10054.375
    "H123" = 23056.682
    "J123" = 23056.682
    "L123" = 23056.682
    "N123" = -27956.682
    "H132" = 2450.5715
    "I132" = 2032.9048
    "J132" = 3077.0715
    "K132" = 6098.7144
    "L132" = 9231.2145
    "M132" = -3568.7144
    "N132" = -14291.2145
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H63" = 21000.6
    "J63" = 21000.6
    "L63" = 21000.6
    "N63" = -22498.6
    "H66" = 21000.6
    "J66" = 21000.6
    "L66" = 63001.8
    "N66" = -70489.79999999999
    "H122" = 2471658.8
    "I122" = 2980955.8
    "K122" = 8942867.399999999
    "M122" = -8940417.399999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H62" = 3500
    "I62" = 3500
    "K62" = 3500
    "M62" = -2876
    "H65" = 3500
    "I65" = 3500
    "K65" = 17500
    "M65" = -14380
    "H132" = 1273.6177
    "I132" = 970.6667
    "J132" = 1763
    "K132" = 2912.0001
    "L132" = 5289
    "M132" = -382.0001000000002
    "N132" = -10349
    "H136" = 8199117
    "I136" = 2618.195
    "J136" = 25001940
    "K136" = 7854.585000000001
    "L136" = 75005820
    "M136" = -5304.585000000001
    "N136" = -75010920
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
